{"js": "const replacements = [\n    [\"2023-08-07 Monday\", \"2023-08-08 Tuesday\"],\n    [\"31\u00d738=\", \"90\u00d795=\"],\n    [\"90\u00d774=\", \"24\u00d722=\"],\n    [\"36\u00d718=\", \"72\u00d794=\"],\n    [\"49\u00d738=\", \"58\u00d714=\"],\n    [\"47\u00d756=\", \"62\u00d748=\"],\n    [\"44\u00d768=\", \"13\u00d720=\"],\n    [\"15\u00d725=\", \"77\u00d771=\"],\n    [\"42\u00d737=\", \"91\u00d745=\"],\n    [\"55\u00d737=\", \"97\u00d758=\"],\n    [\"25\u00d768=\", \"92\u00d799=\"],\n    [\"29\u00d748=\", \"51\u00d786=\"],\n    [\"98\u00d713=\", \"29\u00d759=\"],\n    [\"60\u00d794=\", \"14\u00d763=\"],\n    [\"60\u00d718=\", \"89\u00d758=\"],\n    [\"61\u00d717=\", \"32\u00d721=\"],\n    [\"85\u00d752=\", \"96\u00d749=\"],\n    [\"98\u00d789=\", \"62\u00d757=\"],\n    [\"16\u00d745=\", \"42\u00d731=\"],\n    [\"39\u00d750=\", \"19\u00d729=\"],\n    [\"90\u00d783=\", \"49\u00d781=\"],\n    [\"40\u00d759=\", \"25\u00d765=\"],\n    [\"68\u00d732=\", \"47\u00d712=\"],\n    [\"85\u00d767=\", \"43\u00d732=\"],\n    [\"47\u00d774=\", \"11\u00d760=\"],\n    [\"38\u00d767=\", \"79\u00d736=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Map of old text -> new text (date line + multiplication problems in the table)\n$replacements = @{\n    \"2023-08-07 Monday\" = \"2023-08-08 Tuesday\"\n    \"31\u00d738=\"            = \"90\u00d795=\"\n    \"90\u00d774=\"            = \"24\u00d722=\"\n    \"36\u00d718=\"            = \"72\u00d794=\"\n    \"49\u00d738=\"            = \"58\u00d714=\"\n    \"47\u00d756=\"            = \"62\u00d748=\"\n    \"44\u00d768=\"            = \"13\u00d720=\"\n    \"15\u00d725=\"            = \"77\u00d771=\"\n    \"42\u00d737=\"            = \"91\u00d745=\"\n    \"55\u00d737=\"            = \"97\u00d758=\"\n    \"25\u00d768=\"            = \"92\u00d799=\"\n    \"29\u00d748=\"            = \"51\u00d786=\"\n    \"98\u00d713=\"            = \"29\u00d759=\"\n    \"60\u00d794=\"            = \"14\u00d763=\"\n    \"60\u00d718=\"            = \"89\u00d758=\"\n    \"61\u00d717=\"            = \"32\u00d721=\"\n    \"85\u00d752=\"            = \"96\u00d749=\"\n    \"98\u00d789=\"            = \"62\u00d757=\"\n    \"16\u00d745=\"            = \"42\u00d731=\"\n    \"39\u00d750=\"            = \"19\u00d729=\"\n    \"90\u00d783=\"            = \"49\u00d781=\"\n    \"40\u00d759=\"            = \"25\u00d765=\"\n    \"68\u00d732=\"            = \"47\u00d712=\"\n    \"85\u00d767=\"            = \"43\u00d732=\"\n    \"47\u00d774=\"            = \"11\u00d760=\"\n    \"38\u00d767=\"            = \"79\u00d736=\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n    $rng = $d.Content\n    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
